$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1/J1 - copy formatting (bold/border/centered) from existing header H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-17 for new columns I (I0) and J (IF)
$data = @(
    @(1, 4),
    @(6, 9),
    @(1, 3),
    @(1, 4),
    @(1, 6),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(1, 6),
    @(1, 4),
    @(1, 3),
    @(1, 4),
    @(3, 5),
    @(2, 4),
    @(1, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
